$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.511.06"
$ws.Range("E2").Value = "  -6.85%  "
$ws.Range("D3").Value = "2.587.72"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'301.03"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("D6").Value = "'96.24"
$ws.Range("E6").Value = "  -5.10%  "
$ws.Range("E7").Value = "  -4.32%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").Value = "'36.72"
$ws.Range("E10").Value = "  -7.70%  "
$ws.Range("D11").Value = "'0.0817"
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("D12").Value = "'7.79"
$ws.Range("E12").Value = "  -5.18%  "
$ws.Range("D13").Value = "2.986.82"
$ws.Range("E13").Value = "  -0.67%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "2.590.12"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").Value = "'0.890"
$ws.Range("E16").Value = "  -4.37%  "
$ws.Range("E17").Value = "  -4.99%  "
$ws.Range("D18").Value = "43.508.86"
$ws.Range("E18").Value = "  -7.00%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0980"
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.66"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'12.33"
$ws.Range("E21").Value = "  -5.98%  "
$ws.Range("D22").Value = "'72.83"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "'264.96"
$ws.Range("E23").Value = "  -5.10%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.93"
$ws.Range("E24").Value = "  -4.13%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "'29.20"
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'10.28"
$ws.Range("E28").Value = "  -4.06%  "
$ws.Range("D29").Value = "'37.96"
$ws.Range("E29").Value = "  -3.65%  "
$ws.Range("E30").Value = "  -6.76%  "
$ws.Range("E31").Value = "  -5.97%  "
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").Value = "'152.44"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'2.79"
$ws.Range("E35").Value = "  -2.18%  "
$ws.Range("D36").Value = "'0.0810"
$ws.Range("E36").Value = "  -4.45%  "
$ws.Range("E37").Value = "  -4.73%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.121"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'24.33"
$ws.Range("E39").Value = "  +4.43%  "
$ws.Range("D40").Value = "'16.62"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").Value = "'3.61"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("E42").Value = "  -5.75%  "
$ws.Range("D43").Value = "'3.87"
$ws.Range("E43").Value = "  -5.86%  "
$ws.Range("D44").Value = "2.043.36"
$ws.Range("E44").Value = "  -3.72%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'87.91"
$ws.Range("E46").Value = "  -6.21%  "
$ws.Range("D47").Value = "'9.08"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "2.843.90"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'1.60"
$ws.Range("E49").Value = "  +2.86%  "
$ws.Range("D50").Value = "'105.59"
$ws.Range("E50").Value = "  -3.87%  "
$ws.Range("E51").Value = "  -5.62%  "
